# Insert a new data row at row 43 (shifts existing rows 43:100 down to 44:101)
# and populate it with a new Membrillo price record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("43:43").Insert()

$ws.Range("A43").Value = 8
$ws.Range("B43").Value = "Terminal La Palmera de La Serena"
$ws.Range("C43").Value = "Coquimbo"
$ws.Range("D43").Value = 45082
$ws.Range("E43").Value = 4
$ws.Range("F43").Value = "Fruta"
$ws.Range("G43").Value = 100104
$ws.Range("H43").Value = "Frutos de pepita"
$ws.Range("I43").Value = 100104003
$ws.Range("J43").Value = "Membrillo"
$ws.Range("K43").Value = "Champion"
$ws.Range("L43").Value = "Primera"
$ws.Range("M43").Value = 10
$ws.Range("N43").Value = 200000
$ws.Range("O43").Value = 210000
$ws.Range("P43").Value = 205000
$ws.Range("Q43").Value = "`$/bins (450 kilos)"
$ws.Range("R43").Value = "Región de O'Higgins"
$ws.Range("S43").Value = 456
$ws.Range("T43").Value = 450
